# Update Sheets: refresh market-price-derived columns (H:N) for the specific
# leve rows called out in the scheduled-runner diff. Values are plain numbers
# (no formulas in this workbook), so we just overwrite the stored cell values.
# A handful of rows additionally gain or lose a cell entirely (e.g. LeveProfitHQ
# becomes applicable/inapplicable), handled with .Value assignment / ClearContents().
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 31
$ws.Range("H31").Value = 16694.5
$ws.Range("I31").Value = 16694.5
$ws.Range("K31").Value = 50083.5
$ws.Range("M31").Value = -49853.5

# Row 113
$ws.Range("H113").Value = 3500.625
$ws.Range("I113").Value = 2502.5
$ws.Range("J113").Value = 3833.3333
$ws.Range("K113").Value = 2502.5
$ws.Range("L113").Value = 3833.3333
$ws.Range("M113").Value = 751.5
$ws.Range("N113").Value = -10341.3333

# Row 125
$ws.Range("H125").Value = 60290.777
$ws.Range("I125").Value = 128207.75
$ws.Range("K125").Value = 1153869.75
$ws.Range("M125").Value = -1151409.75

# Row 132
$ws.Range("H132").Value = 1976.5333
$ws.Range("I132").Value = 1524.5385
$ws.Range("J132").Value = 2595.0527
$ws.Range("K132").Value = 4573.6155
$ws.Range("L132").Value = 7785.158100000001
$ws.Range("M132").Value = -2043.6155
$ws.Range("N132").Value = -12845.1581

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 2038
$ws.Range("I2").Value = 1547.625
$ws.Range("J2").Value = 2822.6
$ws.Range("K2").Value = 1547.625
$ws.Range("L2").Value = 2822.6
$ws.Range("M2").Value = -1434.625
$ws.Range("N2").Value = -3048.6

# Row 4
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("M4").ClearContents()

# Row 5
$ws.Range("H5").Value = 480
$ws.Range("I5").Value = 480
$ws.Range("K5").Value = 480
$ws.Range("M5").Value = -368

# Row 6
$ws.Range("H6").Value = 9883.333000000001
$ws.Range("I6").Value = 9883.333000000001
$ws.Range("K6").Value = 9883.333000000001
$ws.Range("M6").Value = -9710.333000000001

# Row 44
$ws.Range("H44").Value = 32666.666
$ws.Range("J44").Value = 39000
$ws.Range("L44").Value = 39000
$ws.Range("N44").Value = -39976

# Row 74
$ws.Range("H74").Value = 1068.8636
$ws.Range("I74").Value = 677.35297
$ws.Range("J74").Value = 2400
$ws.Range("K74").Value = 677.35297
$ws.Range("L74").Value = 2400
$ws.Range("M74").Value = 196.64703
$ws.Range("N74").Value = -4148

# Row 77
$ws.Range("H77").Value = 1068.8636
$ws.Range("I77").Value = 677.35297
$ws.Range("J77").Value = 2400
$ws.Range("K77").Value = 3386.76485
$ws.Range("L77").Value = 12000
$ws.Range("M77").Value = 981.23515
$ws.Range("N77").Value = -20736

# Row 110
$ws.Range("H110").Value = 745.94116
$ws.Range("I110").Value = 573.75
$ws.Range("J110").Value = 1159.2
$ws.Range("K110").Value = 573.75
$ws.Range("L110").Value = 1159.2
$ws.Range("M110").Value = 1471.25
$ws.Range("N110").Value = -5249.2

# Row 116
$ws.Range("H116").Value = 2038
$ws.Range("I116").Value = 1547.625
$ws.Range("J116").Value = 2822.6
$ws.Range("K116").Value = 1547.625
$ws.Range("L116").Value = 2822.6
$ws.Range("M116").Value = 746.375
$ws.Range("N116").Value = -7410.6

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 2038
$ws.Range("I3").Value = 1547.625
$ws.Range("J3").Value = 2822.6
$ws.Range("K3").Value = 1547.625
$ws.Range("L3").Value = 2822.6
$ws.Range("M3").Value = -1433.625
$ws.Range("N3").Value = -3050.6

# Row 4
$ws.Range("H4").Value = 480
$ws.Range("I4").Value = 480
$ws.Range("K4").Value = 480
$ws.Range("M4").Value = -365

# Row 22
$ws.Range("H22").Value = 12875
$ws.Range("I22").Value = 12875
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 12875
$ws.Range("L22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("N22").Value = -12702

# Row 134
$ws.Range("H134").Value = 2105.6765
$ws.Range("I134").Value = 1987.2812
$ws.Range("J134").Value = 4000
$ws.Range("K134").Value = 5961.8436
$ws.Range("L134").Value = 12000
$ws.Range("M134").Value = -3426.8436
$ws.Range("N134").Value = -17070

$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 492.76923
$ws.Range("I22").Value = 261.4
$ws.Range("J22").Value = 1264
$ws.Range("K22").Value = 261.4
$ws.Range("L22").Value = 1264
$ws.Range("M22").Value = 88.60000000000002
$ws.Range("N22").Value = -1964

# Row 23
$ws.Range("H23").Value = 50000
$ws.Range("I23").Value = 50000
$ws.Range("K23").Value = 50000
$ws.Range("M23").Value = -49760

# Row 27
$ws.Range("H27").Value = 50000
$ws.Range("I27").Value = 50000
$ws.Range("K27").Value = 50000
$ws.Range("M27").Value = -49808

# Row 31
$ws.Range("H31").Value = 2287.4
$ws.Range("I31").Value = 1627.8572
$ws.Range("K31").Value = 1627.8572
$ws.Range("M31").Value = -1332.8572

# Row 34
$ws.Range("H34").Value = 2287.4
$ws.Range("I34").Value = 1627.8572
$ws.Range("K34").Value = 1627.8572
$ws.Range("M34").Value = -1425.8572

# Row 99
$ws.Range("H99").Value = 2199
$ws.Range("I99").Value = 2000
$ws.Range("J99").Value = 2331.6667
$ws.Range("K99").Value = 2000
$ws.Range("L99").Value = 2331.6667
$ws.Range("M99").Value = -502
$ws.Range("N99").Value = -5327.6667

# Row 126
$ws.Range("H126").Value = 2199
$ws.Range("I126").Value = 2000
$ws.Range("J126").Value = 2331.6667
$ws.Range("K126").Value = 6000
$ws.Range("L126").Value = 6995.000100000001
$ws.Range("M126").Value = -3530
$ws.Range("N126").Value = -11935.0001

$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 4832.222
$ws.Range("I4").Value = 550
$ws.Range("J4").Value = 6055.7144
$ws.Range("K4").Value = 1650
$ws.Range("L4").Value = 18167.1432
$ws.Range("M4").Value = -1538
$ws.Range("N4").Value = -18391.1432

$ws = $wb.Worksheets.Item("GSM")
# Row 29
$ws.Range("H29").Value = 300000
$ws.Range("I29").Value = 300000
$ws.Range("K29").Value = 300000
$ws.Range("M29").Value = -299710

$ws = $wb.Worksheets.Item("LTW")
# Row 9
$ws.Range("H9").Value = 3331.3333
$ws.Range("I9").Value = 481.33334
$ws.Range("K9").Value = 481.33334
$ws.Range("M9").Value = -257.33334

# Row 16
$ws.Range("H16").Value = 776.86206
$ws.Range("I16").Value = 771.37036
$ws.Range("J16").Value = 851
$ws.Range("K16").Value = 771.37036
$ws.Range("L16").Value = 851
$ws.Range("M16").Value = -601.37036
$ws.Range("N16").Value = -1191

# Row 128
$ws.Range("H128").Value = 265429
$ws.Range("J128").Value = 265429
$ws.Range("L128").Value = 265429
$ws.Range("N128").Value = -275389

# Row 132
$ws.Range("H132").Value = 5217.4287
$ws.Range("I132").Value = 5910.5884
$ws.Range("J132").Value = 4146.1816
$ws.Range("K132").Value = 17731.7652
$ws.Range("L132").Value = 12438.5448
$ws.Range("M132").Value = -15201.7652
$ws.Range("N132").Value = -17498.5448

# Row 137
$ws.Range("H137").Value = 29195.6
$ws.Range("J137").Value = 31328.445
$ws.Range("L137").Value = 31328.445
$ws.Range("N137").Value = -41528.445

$ws = $wb.Worksheets.Item("WVR")
# Row 5
$ws.Range("H5").Value = 6535000.5
$ws.Range("I5").Value = 13000000
$ws.Range("J5").Value = 70001
$ws.Range("K5").Value = 13000000
$ws.Range("L5").Value = 70001
$ws.Range("M5").Value = -12999888
$ws.Range("N5").Value = -70225

# Row 21
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("M21").ClearContents()

# Row 29
$ws.Range("H29").Value = 9500
$ws.Range("I29").Value = 4000
$ws.Range("J29").Value = 15000
$ws.Range("K29").Value = 4000
$ws.Range("L29").Value = 15000
$ws.Range("M29").Value = -3710
$ws.Range("N29").Value = -15580

# Row 35
$ws.Range("H35").Value = 0
$ws.Range("I35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("M35").ClearContents()
